$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data rows (4-13) ----
# Values are written in the same order the shared-string table records them
# in the target workbook (reflecting the author's actual editing / revision
# sequence), so the rebuilt xl/sharedStrings.xml lines up index-for-index.
$ws.Range("A4").Value = 4

$ws.Range("C4").Value  = "Dans le sous-dossier /utilisateurs de Memolab se trouve le nouveau fichier à votre nom. Il est désormais accessible par l'usage d'excel ou un éditeur compatible.Vos performances ainsi que votre historique d'activités peuvent y être consultés et partagés."
$ws.Range("B5").Value  = "Menu - Editer les données d'un utilisateur"
$ws.Range("C5").Value  = "Dans le sous-dossier /utilisateurs de Memolab se trouve le nouveau fichier à votre nom. Il est accessible par l'usage d'excel ou un éditeur compatible.Vos performances ainsi que votre historique d'activités peuvent y être consultés, édités et partagés."
$ws.Range("B6").Value  = "Menu - Utilisateur statistiques"
$ws.Range("B4").Value  = "Menu-Création d'un nouvel utilisateur"
$ws.Range("C6").Value  = "Dans le sous-dossier /utilisateurs de Memolab se trouve le fichier à votre nom. Il est accessible par l'usage d'excel ou un éditeur compatible. Vos performances ainsi que votre historique d'activités peuvent y être consultés et partagés."
$ws.Range("B7").Value  = "Menu - Utilisateur partager"
$ws.Range("C7").Value  = "Dans le sous-dossier /utilisateurs de Memolab se trouve le fichier à votre nom. Il est accessible par l'usage d'excel ou un éditeur compatible. Vos performances ainsi que votre historique d'activités peuvent y être consultés et partagés."
$ws.Range("B8").Value  = "Menu - Création fichier de lecon"
$ws.Range("C8").Value  = "Dans le sous-dossier /lecons de Memolab se trouve le fichier lecons-modele.xlsx. Il est accessible par l'usage d'excel ou un éditeur compatible. Il convient de l'ouvrir et de documenter les colonnes des questions et des réponses. Ensuite il faut l'enregistrer sous un nom qui évoque son contenu. Il vous sera loisible de le choisir pour l'étudier dans le cadre des activités de laboratoire."
$ws.Range("B9").Value  = "Menu - Ouvrir lecon"
$ws.Range("C9").Value  = "Dans le sous-dossier /lecons de Memolab se trouvent tous les fichiers des leçons. Ils sont accessibles par l'usage d'excel ou un éditeur compatible. Ils peuvent être consultés, modifiés et enregistrés."
$ws.Range("B10").Value = "Menu - Supprimer une lecon"
$ws.Range("C10").Value = "Dans le sous-dossier /lecons de Memolab se trouvent tous les fichiers des lecons. Il sont accessibles depuis l'explorateur de fichiers de windows / Finder. Ils peuvent y être supprimés au moyen des commandes du système."
$ws.Range("B11").Value = "Menu - Partager une lecon"
$ws.Range("C11").Value = "Dans le sous-dossier /lecons de Memolab se trouvent tous les fichiers des lecons. Il sont accessibles depuis l'explorateur de fichiers de windows / Finder. Ils peuvent y être copiés et partagés au moyen des commandes du système."
$ws.Range("C12").Value = "Les paramètres généraux des préférences de Memolab sont tous dans le fichiers excel nommée params_generaux.xlsx qui se trouve dans le répertoire des scripts de l'application. Ils peuvent êtres édités à partir de excel ou tout autre éditeur compatible. A modifier avec grande prudence !"
$ws.Range("B12").Value = "Menu - Options préférences"
$ws.Range("B13").Value = "Menu - A propos de"
$ws.Range("C13").Value = " Memolab 2021 version a * Auteur : Jean-Claude Vouillamoz * mail : jcvouillamoz@gmail.com * Phone : + 41 79 212 84 52 * License GNU"

# ---- Row heights (wrapped text autosize, matching the saved heights) ----
$ws.Rows.Item(4).RowHeight  = 43.2
$ws.Rows.Item(5).RowHeight  = 43.2
$ws.Rows.Item(6).RowHeight  = 43.2
$ws.Rows.Item(7).RowHeight  = 43.2
$ws.Rows.Item(8).RowHeight  = 72
$ws.Rows.Item(9).RowHeight  = 43.2
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 57.6
$ws.Rows.Item(13).RowHeight = 28.8

# ---- Column B width (grew to fit the new, wider menu labels) ----
$ws.Columns.Item(2).ColumnWidth = 34

# ---- Selection / scroll position left where editing ended ----
$ws.Range("C13").Select() | Out-Null
